$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Drop the standalone "Meta description: ..." paragraph that used
#    to sit right under the H1 title.
# ------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Meta description")) {
        $p.Range.Delete()
    }
}

# ------------------------------------------------------------------
# 2) Replace the trailing "Please create a feature image..." note
#    (the last paragraph in the doc) with two new paragraphs that
#    carry the title (bold) and the meta description text (italic) -
#    i.e. move/reshape the content removed in step 1 down to the end.
# ------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text.StartsWith("Please create a feature image")) {
        $p.Range.Delete()
    }
}

$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$insertPoint = $d.Range($lastPara.Range.End, $lastPara.Range.End)

$titleText = "Play Alkemor's Tower for Free - Exciting Gameplay with Elementals"
$descText  = "Find out about the exciting gameplay and free spins in Alkemor's Tower, set in a magical world with mythological elementals. Play for free now."

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
       '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
       '<pkg:xmlData>' +
       '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
       '<w:body>' +
       '<w:p><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>' + $titleText + '</w:t></w:r></w:p>' +
       '<w:p><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>' + $descText + '</w:t></w:r></w:p>' +
       '</w:body></w:document>' +
       '</pkg:xmlData></pkg:part></pkg:package>'

[void]$insertPoint.InsertXML($xml)
